# Apply updated dSF (column F) values to Sheet1, reflecting repulled/pushed data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = 2
    4  = -1
    6  = -3
    7  = -4
    8  = 2
    9  = 1
    10 = 3
    11 = -4
    13 = 2
    14 = -1
    15 = -3
    16 = -2
    17 = 1
    18 = 1
    19 = -4
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
